{"js": "/*\n * Edit applied (per the provided diff):\n *  1) Merge the run/proofErr fragments around \"lootbox\" into a single run.\n *  2) Merge the run/proofErr fragments around \"Skinport\"/\"Skinbaron\" into a single run.\n *  3) Merge the three runs forming \"...with a \"+\"choice for their \"+\"starting amount...\"\n *     into a single run.\n *  4) Expand the single run \"from the site CS:GO Backpack\" into a longer sentence,\n *     split across many runs, describing the API in more detail.\n *\n * Implementation notes:\n *   Office.js's Range.search()/insertText() merge/clean-up behavior for runs that\n *   span w:proofErr markers is inconsistent in this host (it can reorder sibling\n *   runs when a replaced range crosses more than one existing run). To apply the\n *   edit reliably and deterministically we:\n *     - locate each target paragraph by its stable w14:paraId (Office.js:\n *       `paragraph.uniqueLocalId`),\n *     - rebuild that paragraph's full, exact OOXML (preserving every untouched\n *       sibling run and its rsidR/other attributes),\n *     - replace the whole paragraph range in one shot via `range.insertOoxml(...,\n *       Word.InsertLocation.replace)`, which is a supported Word JS API call and\n *       (unlike sub-paragraph range replace here) does not reorder content.\n */\n\nfunction flatOpcPackage(paragraphXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + paragraphXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\n// Replacement OOXML for each edited paragraph, keyed by its stable w14:paraId.\nconst REPLACEMENTS = {\n  // \"Within the game Counter-Strike... lootbox ... Skinport ... Skinbaron ...\"\n  \"2036E964\": \"<w:p w14:paraId=\\\"2036E964\\\" w14:textId=\\\"334B724A\\\" w:rsidR=\\\"00160876\\\" w:rsidRPr=\\\"00160876\\\" w:rsidRDefault=\\\"003B023A\\\" w:rsidP=\\\"00160876\\\"><w:r><w:t xml:space=\\\"preserve\\\">Within the game Counter-Strike: Global Offensive (CS:GO) there exists a \\u201clootbox\\u201d system where players pay real money for a case and a key to roll for a random item contained within, with the potential of receiving </w:t></w:r><w:r w:rsidR=\\\"000269A9\\\"><w:t>a rare item worth lots of money on</w:t></w:r><w:r w:rsidR=\\\"00E75903\\\"><w:t xml:space=\\\"preserve\\\"> the</w:t></w:r><w:r w:rsidR=\\\"000269A9\\\"><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r w:rsidR=\\\"00E75903\\\"><w:t>in-game</w:t></w:r><w:r w:rsidR=\\\"000269A9\\\"><w:t xml:space=\\\"preserve\\\"> marketplace </w:t></w:r><w:r w:rsidR=\\\"00E75903\\\"><w:t>called the \\u201c</w:t></w:r><w:r w:rsidR=\\\"000269A9\\\"><w:t>Steam Market</w:t></w:r><w:r w:rsidR=\\\"00E75903\\\"><w:t>\\u201d</w:t></w:r><w:r w:rsidR=\\\"000269A9\\\"><w:t xml:space=\\\"preserve\\\"> or external sites such as \\u201cSkinport\\u201d or \\u201cSkinbaron\\u201d </w:t></w:r><w:r w:rsidR=\\\"009C387B\\\"><w:t xml:space=\\\"preserve\\\">which allow for users to cash out their items for real money by selling them to other players or the marketplace itself. </w:t></w:r><w:r w:rsidR=\\\"005F5B2E\\\"><w:t>This system presents the curiosity, if one were to have a large sum of money, how profitable would gambling it on CS:GO be? Which choice of case is most profitable</w:t></w:r><w:r w:rsidR=\\\"005C6D27\\\"><w:t xml:space=\\\"preserve\\\"> long term</w:t></w:r><w:r w:rsidR=\\\"005F5B2E\\\"><w:t>?</w:t></w:r><w:r w:rsidR=\\\"005C6D27\\\"><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r w:rsidR=\\\"005F5B2E\\\"><w:t xml:space=\\\"preserve\\\">However, there is no capability within the game to simulate such an experience without the risk </w:t></w:r><w:r w:rsidR=\\\"000C1113\\\"><w:t>of spending copious amounts of your own money.</w:t></w:r><w:r w:rsidR=\\\"007D14FE\\\"><w:t xml:space=\\\"preserve\\\"> My solution </w:t></w:r><w:r w:rsidR=\\\"00EC74FC\\\"><w:t>was</w:t></w:r><w:r w:rsidR=\\\"007D14FE\\\"><w:t xml:space=\\\"preserve\\\"> to create a price accurate simulation of the </w:t></w:r><w:r w:rsidR=\\\"00EC74FC\\\"><w:t>case opening system within CS:GO</w:t></w:r><w:r w:rsidR=\\\"006A6B9B\\\"><w:t>, allowing users to purchase cases and sell items much alike to the game</w:t></w:r><w:r w:rsidR=\\\"00CF2BCF\\\"><w:t>, without the aspect of spending actual money.</w:t></w:r></w:p>\",\n  // \"The idea comprises a website... choice for their starting amount of money...\"\n  \"081D9A2B\": \"<w:p w14:paraId=\\\"081D9A2B\\\" w14:textId=\\\"0876F331\\\" w:rsidR=\\\"00E75903\\\" w:rsidRPr=\\\"00E75903\\\" w:rsidRDefault=\\\"00E75903\\\" w:rsidP=\\\"00E75903\\\"><w:r><w:t xml:space=\\\"preserve\\\">The idea comprises a website that simulates the case opening and item system in the game Counter-Strike: Global Offensive (CS:GO). The simulation would provide the user with a choice for their starting amount of money and allow them to purchase cases much like within the game. A case has a range of potential items contained within with varying value attached to each item, when a case is rolled a single item from within its contents is received. The odds of receiving each item in the case will be provided to the user, </w:t></w:r><w:r w:rsidR=\\\"002A6F94\\\"><w:t>higher rarity items such as \\u201cSouvenir Weapons\\u201d, \\u201cKnives\\u201d or \\u201cGloves\\u201d will have much higher value than other items</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">. The simulation will allow the user to sell the items they receive to increase their balance and open more cases. Each case will have a different price depending on the rarity of the items within and the age of the case itself, as well as an accurate average market price for each item which will be scraped from an active source each time the site is loaded. </w:t></w:r><w:r w:rsidR=\\\"00330C99\\\"><w:t>The system will track the users wins and losses on case openings and allow them to look at their statistics once they either choose to end the simulation or run themselves completely out of money.</w:t></w:r></w:p>\",\n  // \"...The aim is to utilise an API from the site CS:GO Backpack\"\n  \"3FACFE7F\": \"<w:p w14:paraId=\\\"3FACFE7F\\\" w14:textId=\\\"5512F8F4\\\" w:rsidR=\\\"003F13EB\\\" w:rsidRPr=\\\"003F13EB\\\" w:rsidRDefault=\\\"003F13EB\\\" w:rsidP=\\\"003F13EB\\\"><w:r><w:t xml:space=\\\"preserve\\\">The main feature objective of the site is to provide a price accurate simulation, with the value of cases and items pegged to the current average prices </w:t></w:r><w:r w:rsidR=\\\"004F0DB6\\\"><w:t xml:space=\\\"preserve\\\">for the item from online marketplaces. The aim is to utilise an API </w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t xml:space=\\\"preserve\\\">from the site </w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t>\\u201c</w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t>CS:GO Backpack</w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t xml:space=\\\"preserve\\\">\\u201d which provides </w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t>average</w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t xml:space=\\\"preserve\\\"> price data</w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t>as JSON objects</w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t xml:space=\\\"preserve\\\"> for</w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t xml:space=\\\"preserve\\\"> all market items</w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t xml:space=\\\"preserve\\\"> from the last 24 hours, week</w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t>and</w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t xml:space=\\\"preserve\\\"> all time, updated every 8 hours</w:t></w:r><w:r w:rsidR=\\\"00C91635\\\"><w:t xml:space=\\\"preserve\\\"> </w:t></w:r></w:p>\"\n};\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"uniqueLocalId\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const replacementXml = REPLACEMENTS[paragraph.uniqueLocalId];\n  if (!replacementXml) continue;\n  const range = paragraph.getRange();\n  range.insertOoxml(flatOpcPackage(replacementXml), Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Edit applied (per the provided diff):\n#  1) Merge the run/proofErr fragments around \"lootbox\" into a single run.\n#  2) Merge the run/proofErr fragments around \"Skinport\"/\"Skinbaron\" into a single run.\n#  3) Merge the three runs forming \"...with a \" + \"choice for their \" + \"starting amount...\"\n#     into a single run.\n#  4) Expand the single run \"from the site CS:GO Backpack\" into a longer sentence,\n#     split across many runs, describing the API in more detail.\n#\n# Implementation notes:\n#   Word's Find/Replace (and Range.Text assignment) across a run that is wrapped by\n#   w:proofErr markers does not reliably drop the now-orphaned w:proofErr elements in\n#   this host, and sub-paragraph Range.InsertXML() calls that span more than one\n#   existing run can reorder sibling runs. To apply the edit deterministically we:\n#     - find each target paragraph using Find (a stable text anchor unique to it),\n#     - take that paragraph's whole Range,\n#     - replace it in one shot via Range.InsertXML(...) with the paragraph's full,\n#       exact OOXML (every untouched sibling run/rsidR preserved) - InsertXML is the\n#       COM analogue of Word JS's Range.insertOoxml and does not reorder content when\n#       applied to a whole paragraph range.\n\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphOoxml([string]$anchorText, [string]$ooxml) {\n    $find = $d.Content\n    $find.Find.ClearFormatting()\n    $find.Find.Execute($anchorText, $true) | Out-Null\n    if (-not $find.Find.Found) {\n        throw \"Anchor text not found: $anchorText\"\n    }\n    $targetParagraph = $find.Paragraphs(1)\n    $targetParagraph.Range.InsertXML($ooxml)\n}\n\n# Paragraph: \"Within the game Counter-Strike... lootbox ... Skinport ... Skinbaron ...\"\n$xmlLootbox = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p w14:paraId=\"2036E964\" w14:textId=\"334B724A\" w:rsidR=\"00160876\" w:rsidRPr=\"00160876\" w:rsidRDefault=\"003B023A\" w:rsidP=\"00160876\"><w:r><w:t xml:space=\"preserve\">Within the game Counter-Strike: Global Offensive (CS:GO) there exists a \u201clootbox\u201d system where players pay real money for a case and a key to roll for a random item contained within, with the potential of receiving </w:t></w:r><w:r w:rsidR=\"000269A9\"><w:t>a rare item worth lots of money on</w:t></w:r><w:r w:rsidR=\"00E75903\"><w:t xml:space=\"preserve\"> the</w:t></w:r><w:r w:rsidR=\"000269A9\"><w:t xml:space=\"preserve\"> </w:t></w:r><w:r w:rsidR=\"00E75903\"><w:t>in-game</w:t></w:r><w:r w:rsidR=\"000269A9\"><w:t xml:space=\"preserve\"> marketplace </w:t></w:r><w:r w:rsidR=\"00E75903\"><w:t>called the \u201c</w:t></w:r><w:r w:rsidR=\"000269A9\"><w:t>Steam Market</w:t></w:r><w:r w:rsidR=\"00E75903\"><w:t>\u201d</w:t></w:r><w:r w:rsidR=\"000269A9\"><w:t xml:space=\"preserve\"> or external sites such as \u201cSkinport\u201d or \u201cSkinbaron\u201d </w:t></w:r><w:r w:rsidR=\"009C387B\"><w:t xml:space=\"preserve\">which allow for users to cash out their items for real money by selling them to other players or the marketplace itself. </w:t></w:r><w:r w:rsidR=\"005F5B2E\"><w:t>This system presents the curiosity, if one were to have a large sum of money, how profitable would gambling it on CS:GO be? Which choice of case is most profitable</w:t></w:r><w:r w:rsidR=\"005C6D27\"><w:t xml:space=\"preserve\"> long term</w:t></w:r><w:r w:rsidR=\"005F5B2E\"><w:t>?</w:t></w:r><w:r w:rsidR=\"005C6D27\"><w:t xml:space=\"preserve\"> </w:t></w:r><w:r w:rsidR=\"005F5B2E\"><w:t xml:space=\"preserve\">However, there is no capability within the game to simulate such an experience without the risk </w:t></w:r><w:r w:rsidR=\"000C1113\"><w:t>of spending copious amounts of your own money.</w:t></w:r><w:r w:rsidR=\"007D14FE\"><w:t xml:space=\"preserve\"> My solution </w:t></w:r><w:r w:rsidR=\"00EC74FC\"><w:t>was</w:t></w:r><w:r w:rsidR=\"007D14FE\"><w:t xml:space=\"preserve\"> to create a price accurate simulation of the </w:t></w:r><w:r w:rsidR=\"00EC74FC\"><w:t>case opening system within CS:GO</w:t></w:r><w:r w:rsidR=\"006A6B9B\"><w:t>, allowing users to purchase cases and sell items much alike to the game</w:t></w:r><w:r w:rsidR=\"00CF2BCF\"><w:t>, without the aspect of spending actual money.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\nSet-ParagraphOoxml \"Within the game Counter-Strike: Global Offensive (CS:GO) there exists a\" $xmlLootbox\n\n# Paragraph: \"The idea comprises a website... choice for their starting amount of money...\"\n$xmlChoice = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p w14:paraId=\"081D9A2B\" w14:textId=\"0876F331\" w:rsidR=\"00E75903\" w:rsidRPr=\"00E75903\" w:rsidRDefault=\"00E75903\" w:rsidP=\"00E75903\"><w:r><w:t xml:space=\"preserve\">The idea comprises a website that simulates the case opening and item system in the game Counter-Strike: Global Offensive (CS:GO). The simulation would provide the user with a choice for their starting amount of money and allow them to purchase cases much like within the game. A case has a range of potential items contained within with varying value attached to each item, when a case is rolled a single item from within its contents is received. The odds of receiving each item in the case will be provided to the user, </w:t></w:r><w:r w:rsidR=\"002A6F94\"><w:t>higher rarity items such as \u201cSouvenir Weapons\u201d, \u201cKnives\u201d or \u201cGloves\u201d will have much higher value than other items</w:t></w:r><w:r><w:t xml:space=\"preserve\">. The simulation will allow the user to sell the items they receive to increase their balance and open more cases. Each case will have a different price depending on the rarity of the items within and the age of the case itself, as well as an accurate average market price for each item which will be scraped from an active source each time the site is loaded. </w:t></w:r><w:r w:rsidR=\"00330C99\"><w:t>The system will track the users wins and losses on case openings and allow them to look at their statistics once they either choose to end the simulation or run themselves completely out of money.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\nSet-ParagraphOoxml \"The idea comprises a website\" $xmlChoice\n\n# Paragraph: \"...The aim is to utilise an API from the site CS:GO Backpack\"\n$xmlBackpack = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p w14:paraId=\"3FACFE7F\" w14:textId=\"5512F8F4\" w:rsidR=\"003F13EB\" w:rsidRPr=\"003F13EB\" w:rsidRDefault=\"003F13EB\" w:rsidP=\"003F13EB\"><w:r><w:t xml:space=\"preserve\">The main feature objective of the site is to provide a price accurate simulation, with the value of cases and items pegged to the current average prices </w:t></w:r><w:r w:rsidR=\"004F0DB6\"><w:t xml:space=\"preserve\">for the item from online marketplaces. The aim is to utilise an API </w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t xml:space=\"preserve\">from the site </w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t>\u201c</w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t>CS:GO Backpack</w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t xml:space=\"preserve\">\u201d which provides </w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t>average</w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t xml:space=\"preserve\"> price data</w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t xml:space=\"preserve\"> </w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t>as JSON objects</w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t xml:space=\"preserve\"> for</w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t xml:space=\"preserve\"> all market items</w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t xml:space=\"preserve\"> from the last 24 hours, week</w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t xml:space=\"preserve\"> </w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t>and</w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t xml:space=\"preserve\"> all time, updated every 8 hours</w:t></w:r><w:r w:rsidR=\"00C91635\"><w:t xml:space=\"preserve\"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\nSet-ParagraphOoxml \"The main feature objective of the site\" $xmlBackpack\n"}
